$d = $word.ActiveDocument

# Identify the two list-bullet paragraphs that mention the inline
# affiliate disclosure requirement; these are being removed because the
# disclosure now lives in the website sidebar instead of inline copy.
$targets = @(
    "✅ Introduction with affiliate disclosure",
    "Affiliate Disclosure: Required in introduction paragraph"
)

$toDelete = @()
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r", "`a")
    foreach ($target in $targets) {
        if ($t -eq $target) {
            $toDelete += $p
        }
    }
}

# Delete from the bottom of the document upward so earlier deletions
# don't shift the character offsets backing the ranges queued later.
for ($i = $toDelete.Count - 1; $i -ge 0; $i--) {
    $toDelete[$i].Range.Delete()
}
